$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Revenue")
$products = @{
    2 = "Laptop Pro 15"
    3 = "IT Support Service"
    4 = "Printer Paper"
    5 = "Laptop Pro 15"
    6 = "Printer Paper"
    7 = "Office Chair"
    8 = "Laptop Pro 15"
    9 = "Laptop Pro 15"
    10 = "IT Support Service"
    11 = "Printer Paper"
    12 = "Laptop Pro 15"
    13 = "Laptop Pro 15"
    14 = "IT Support Service"
    15 = "Laptop Pro 15"
    16 = "Laptop Pro 15"
    17 = "Laptop Pro 15"
    18 = "Printer Paper"
    19 = "IT Support Service"
    20 = "Printer Paper"
    21 = "Office Chair"
    22 = "Printer Paper"
    23 = "IT Support Service"
    24 = "Laptop Pro 15"
    25 = "Printer Paper"
    26 = "Laptop Pro 15"
    27 = "Printer Paper"
    28 = "Laptop Pro 15"
    29 = "IT Support Service"
    30 = "IT Support Service"
    31 = "IT Support Service"
    32 = "Office Chair"
    33 = "Laptop Pro 15"
    34 = "IT Support Service"
    35 = "Laptop Pro 15"
    36 = "Laptop Pro 15"
    37 = "Marketing Consultation"
    38 = "Office Chair"
    39 = "Marketing Consultation"
    40 = "Office Chair"
    41 = "Laptop Pro 15"
    42 = "Office Chair"
    43 = "Marketing Consultation"
    44 = "Office Chair"
    45 = "Office Chair"
    46 = "Office Chair"
    47 = "Marketing Consultation"
    48 = "Marketing Consultation"
    49 = "IT Support Service"
    50 = "IT Support Service"
    51 = "Marketing Consultation"
}
foreach ($row in $products.Keys) {
    $ws.Cells.Item($row, 4).Value = $products[$row]
}

$ws = $wb.Worksheets.Item("Expenses")
$products = @{
    2 = "Office Supplies Bundle"
    3 = "Office Supplies Bundle"
    4 = "Marketing Consultation"
    5 = "IT Support Service"
    6 = "IT Support Service"
    7 = "Cleaning Supplies Kit"
    8 = "Marketing Consultation"
    9 = "IT Support Service"
    10 = "IT Support Service"
    11 = "Office Supplies Bundle"
    12 = "Office Supplies Bundle"
    13 = "Office Supplies Bundle"
    14 = "Office Supplies Bundle"
    15 = "IT Support Service"
    16 = "IT Support Service"
    17 = "Marketing Consultation"
    18 = "Office Supplies Bundle"
    19 = "Marketing Consultation"
    20 = "Office Supplies Bundle"
    21 = "Office Supplies Bundle"
    22 = "Cleaning Supplies Kit"
    23 = "Cleaning Supplies Kit"
    24 = "Cleaning Supplies Kit"
    25 = "Office Supplies Bundle"
    26 = "Cleaning Supplies Kit"
    27 = "Cleaning Supplies Kit"
    28 = "Marketing Consultation"
    29 = "IT Support Service"
    30 = "Office Supplies Bundle"
    31 = "Cleaning Supplies Kit"
    32 = "IT Support Service"
    33 = "Marketing Consultation"
    34 = "IT Support Service"
    35 = "Marketing Consultation"
    36 = "IT Support Service"
}
foreach ($row in $products.Keys) {
    $ws.Cells.Item($row, 4).Value = $products[$row]
}

